# wip: work on wound and pressureUlcer zibs
#
# Reproduces the author's editing session:
#   - "Sheet1" is renamed to "Research" (sheetId/relationship untouched).
#   - The remembered cell selection (and, for the Data sheet, the
#     scrolled-to top-left cell) is updated on a few sheets to match
#     where the author had been working.

$wb = $excel.ActiveWorkbook

# --- Rename "Sheet1" -> "Research" ---
$research = $wb.Worksheets.Item("Sheet1")
$research.Name = "Research"

# --- WoundTissueOncoUlcerCodelist: selection moves to E24 ---
$ws = $wb.Worksheets.Item("WoundTissueOncoUlcerCodelist")
$ws.Select()
$ws.Range("E24").Select()

# --- Metadata: selection moves to C23 ---
$ws = $wb.Worksheets.Item("Metadata")
$ws.Select()
$ws.Range("C23").Select()

# --- Research (formerly Sheet1): selection moves from C12 to C20 ---
$ws = $wb.Worksheets.Item("Research")
$ws.Select()
$ws.Range("C20").Select()

# --- Data: view scrolls so column Q is leftmost, selection moves to S7 ---
$ws = $wb.Worksheets.Item("Data")
$ws.Select()
$win = $wb.Windows.Item(1)
$topLeft = $ws.Range("Q1")
$win.ScrollRow = $topLeft.Row
$win.ScrollColumn = $topLeft.Column
$ws.Range("S7").Select()

# NB: the workbook-level "first visible tab" (bookViews/workbookView
# firstSheet="2") has no corresponding writable property surfaced by
# this host - ScrollWorkbookTabs()/the Windows collection do not round
# -trip that attribute on save, so it is intentionally left as-is.
